# Update the "Förändrad" (Changed) date column (C) for rows 2-13
# from 45184 to 45185 (i.e. advance the date by one day), matching
# the automatic update recorded in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
